# "Generate Report for Archive"
#
# The localization status report is being regenerated: rows that were
# previously "Ready for handoff" have moved on to "In Translation", and the
# Status columns (which are sized to fit their widest cell) shrink to match
# the new, shorter text.
#
# Affected columns:
#   - Overview sheet: columns E (zh-cn) and F (de-de)
#   - zh-cn sheet:     column C (Status)
#   - de-de sheet:     column C (Status)

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "In Translation"

# --- Overview sheet --------------------------------------------------------
$wsOverview = $wb.Worksheets.Item(1)
$wsOverview.Range("E2").Value = $newStatus
$wsOverview.Range("F2").Value = $newStatus
$wsOverview.Range("E3").Value = $newStatus
$wsOverview.Range("F3").Value = $newStatus

# --- zh-cn sheet -------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item(2)
$wsZhCn.Range("C2").Value = $newStatus
$wsZhCn.Range("C3").Value = $newStatus

# --- de-de sheet -------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item(3)
$wsDeDe.Range("C2").Value = $newStatus
$wsDeDe.Range("C3").Value = $newStatus

# --- Shrink the Status columns to fit the new, shorter text ----------------
# (the report generator re-sizes each Status column to its content)
$wsOverview.Columns.Item(5).ColumnWidth = 12.5
$wsOverview.Columns.Item(6).ColumnWidth = 12.5
$wsZhCn.Columns.Item(3).ColumnWidth = 12.5
$wsDeDe.Columns.Item(3).ColumnWidth = 12.5
